# Updates crypto price ("D") and 1h volume change ("E") columns
# for rows 2-51 on the active worksheet, reflecting the latest
# GitHub Actions scrape of coinranking.com data.
#
# A handful of price cells (D22, D26, D44, D48) are written with a
# leading apostrophe so Excel keeps the significant trailing zero in
# the text instead of silently normalising it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.178.52'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.658.48'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  +0.41%  '
$ws.Range("D5").Value = '217.06'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '0.5196'
$ws.Range("E6").Value = '  -2.28%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '0.2634'
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("D9").Value = '0.06249'
$ws.Range("E9").Value = '  -2.46%  '
$ws.Range("D10").Value = '20.69'
$ws.Range("E10").Value = '  -5.06%  '
$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").Value = '1.661.49'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").Value = '4.401'
$ws.Range("E13").Value = '  -2.46%  '
$ws.Range("D14").Value = '1.886.80'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '0.5414'
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").Value = '0.0₅8086'
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").Value = '64.27'
$ws.Range("E17").Value = '  -2.12%  '
$ws.Range("D18").Value = '26.209.63'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("D20").Value = '4.611'
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").Value = '190.99'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").Value = "'10.00"
$ws.Range("E22").Value = '  -2.79%  '
$ws.Range("D23").Value = '6.026'
$ws.Range("E23").Value = '  -4.89%  '
$ws.Range("D24").Value = '1.008'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '139.48'
$ws.Range("E25").Value = '  -2.06%  '
$ws.Range("D26").Value = "'0.1220"
$ws.Range("E26").Value = '  -4.62%  '
$ws.Range("D27").Value = '7.129'
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '15.97'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("E29").Value = '  -2.74%  '
$ws.Range("D30").Value = '0.05996'
$ws.Range("E30").Value = '  -4.54%  '
$ws.Range("D31").Value = '1.274'
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").Value = '3.576'
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").Value = '3.246'
$ws.Range("E33").Value = '  -5.96%  '
$ws.Range("D34").Value = '1.612'
$ws.Range("E34").Value = '  -4.50%  '
$ws.Range("D35").Value = '0.9599'
$ws.Range("E35").Value = '  -4.82%  '
$ws.Range("D36").Value = '2.423'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").Value = '2.778'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '0.5636'
$ws.Range("E38").Value = '  -9.03%  '
$ws.Range("D39").Value = '5.991'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("D40").Value = '0.01586'
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("D41").Value = '0.8551'
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("D42").Value = '1.004'
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").Value = '1.013.70'
$ws.Range("E43").Value = '  -7.44%  '
$ws.Range("D44").Value = "'100.00"
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '1.801.66'
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").Value = '  +6.57%  '
$ws.Range("D47").Value = '56.64'
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("D48").Value = "'1.010"
$ws.Range("E48").Value = '  +0.88%  '
$ws.Range("D49").Value = '7.964'
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").Value = '0.05181'
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").Value = '1.463'
$ws.Range("E51").Value = '  -1.26%  '
